$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.906.70"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.348.40"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +9.02%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "622.04"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.17"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +7.14%  "

$ws.Range("E8").Value = "  +2.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.343.64"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.791"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.199"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.707.75"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.76"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.29%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000245"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.07%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.982.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +9.18%  "

$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.346.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +9.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.80"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "483.79"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +10.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.83"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.00%  "

$ws.Range("E23").Value = "  +8.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.78%  "

$ws.Range("E25").Value = "  +2.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.56"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.532.22"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +9.37%  "

$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.188"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.244"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.20%  "

$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +13.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.20"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.11"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.78%  "

$ws.Range("E36").Value = "  -2.44%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "510.31"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.03%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.26"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.72%  "

$ws.Range("E39").Value = "  +3.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.80"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.28%  "

$ws.Range("E41").Value = "  +1.87%  "

$ws.Range("E42").Value = "  -1.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.22"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.785"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +17.14%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.05"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.45"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.18%  "

$ws.Range("E50").Value = "  +6.59%  "

$ws.Range("E51").Value = "  +5.89%  "
